$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order for the data rows
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Full target data for rows 2-11 (Gdf1 -> Acvr1 ligand-receptor pairs, recomputed TPM stats
# plus 5 new rows for the "Resolving-Mac" sending cluster).
$rowsData = @(
    @("FAPs", "Gdf1", "Acvr1", "ECs", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 7.004922666666666, 21.014768, 0.134029393318039, 0.1464771065395205, 0.04812381872, 0.43311436848, 0.02695734204101436, 0.02946095154511537),
    @("FAPs", "Gdf1", "Acvr1", "FAPs", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 19.840792, 59.522376, 0.3796257919253833, 0.4148827819958515, 0.13630624104, 1.22675616936, 0.07635416431558342, 0.08344540540186492),
    @("FAPs", "Gdf1", "Acvr1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 6.293636666666667, 18.88091, 0.12041993100245, 0.1316036924905903, 0.04323728390000001, 0.3891355551, 0.02422006985352437, 0.02646946064965762),
    @("FAPs", "Gdf1", "Acvr1", "MuSCs", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 2, 1, 13.324299, 26.648598, 0.2549418168249328, 0.1857460205306503, 0.09153793413, 0.54922760478, 0.05127646694930142, 0.03735911119376898),
    @("FAPs", "Gdf1", "Acvr1", "Resolving-Mac", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 5.800427666666667, 17.401283, 0.110983066929195, 0.1212903984433873, 0.03984893807, 0.35864044263, 0.02232203266690779, 0.02439514703592444),
    @("Resolving-Mac", "Gdf1", "Acvr1", "ECs", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 7.004922666666666, 21.014768, 0.134029393318039, 0.1464771065395205, 0.1911433248053333, 1.720289923248, 0.1070720512770246, 0.1170161549944051),
    @("Resolving-Mac", "Gdf1", "Acvr1", "FAPs", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 19.840792, 59.522376, 0.3796257919253833, 0.4148827819958515, 0.541395691304, 4.872561221736, 0.3032716276097999, 0.3314373765939866),
    @("Resolving-Mac", "Gdf1", "Acvr1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 6.293636666666667, 18.88091, 0.12041993100245, 0.1316036924905903, 0.1717344637233334, 1.54561017351, 0.09619986114892569, 0.1051342318409327),
    @("Resolving-Mac", "Gdf1", "Acvr1", "MuSCs", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 2, 1, 13.324299, 26.648598, 0.2549418168249328, 0.1857460205306503, 0.363580146813, 2.181480880878, 0.2036653498756315, 0.1483869093368813),
    @("Resolving-Mac", "Gdf1", "Acvr1", "Resolving-Mac", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 5.800427666666667, 17.401283, 0.110983066929195, 0.1212903984433873, 0.1582762697403333, 1.424486427663, 0.08866103426228719, 0.09689525140746291),
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $rowVals[$j]
    }
}
